# Update "Vehicle Sections" sheet (Lox Tank row) with new Mass/Length values,
# and make it the active sheet/selection (replacing "Aerodynamic Properties").

$wb = $excel.ActiveWorkbook

$wsVehicle = $wb.Worksheets.Item("Vehicle Sections")
$wsAero = $wb.Worksheets.Item("Aerodynamic Properties")

# Update the Lox Tank row values
$wsVehicle.Range("B6").Value = 4.42
$wsVehicle.Range("C6").Value = 8.8000000000000007

# Move the selection on the Vehicle Sections sheet to C6
$wsVehicle.Activate() | Out-Null
$wsVehicle.Range("C6").Select() | Out-Null

# Aerodynamic Properties keeps its own selection at C3, but is no longer the active sheet
$wsAero.Range("C3").Select() | Out-Null

# Re-activate Vehicle Sections so it is the active/selected tab on save
$wsVehicle.Activate() | Out-Null
